$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "team_logo" column header in D1 (cell already carries the header style s=1)
$ws.Range("D1").Value = "team_logo"

# New "/Team-logos/logo1.png" values for each team row
$ws.Range("D2").Value = "/Team-logos/logo1.png"
$ws.Range("D3").Value = "/Team-logos/logo1.png"

# Match the data-row style (font/format) used by the neighboring C column cells
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
